# Regenerate save_data column "K" (sheet column G) using updated Strike (K)
# values, as produced by recalculating std/mean and writing the new s_vals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values for rows 2..25 (row 26 is unchanged)
$kValues = @(1, 0, 3, 2, 2, 1, 1, 0, 0, 0, 0, 0, 0, 1, 0, 2, 1, 3, 1, 1, 1, 1, 3, 1)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
